# Update column F ("Lama Proses") values for several rows,
# increasing each by 4 days as reflected in the refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(8, 9, 10, 11, 13, 16, 18, 19, 21, 22, 23)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value2 = $cell.Value2 + 4
}
